$wb = $excel.ActiveWorkbook
$foods = $wb.Worksheets.Item("Foods")

# Add a new "calories" column (column E) to the Foods sheet.
$foods.Range("E1").Value = "calories"

# E2 / E3 are blank (empty string) cells, matching the existing blank
# description_title/description_body cells on row 3. Use the classic
# leading-apostrophe trick to force an explicit empty-text cell instead
# of Excel simply leaving the cell empty/uncreated, then strip the
# "quote prefix" formatting it implies so the cell keeps the workbook's
# default style.
$foods.Range("E2").Value = "'"
$foods.Range("E2").ClearFormats()

$foods.Range("E3").Value = "'"
$foods.Range("E3").ClearFormats()
